$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data per commit diff.
# Price cells that look like plain numbers (e.g. "1.00", "222.51") are
# written with a leading apostrophe so Excel stores them as literal text
# (matching the workbook's original inlineStr cells) instead of silently
# converting them to numeric values and dropping trailing zeros.
$ws.Range('D2').Value = '90.736.03'
$ws.Range('E2').Value = '  +3.63%  '
$ws.Range('D3').Value = '3.210.23'
$ws.Range('E3').Value = '  +1.50%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('D5').Value = '''222.51'
$ws.Range('E5').Value = '  +7.35%  '
$ws.Range('D6').Value = '''641.67'
$ws.Range('E6').Value = '  +5.39%  '
$ws.Range('D7').Value = '''0.404'
$ws.Range('E7').Value = '  +6.94%  '
$ws.Range('D8').Value = '''0.708'
$ws.Range('E8').Value = '  +6.48%  '
$ws.Range('D9').Value = '''0.999'
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('D10').Value = '3.208.41'
$ws.Range('E10').Value = '  +1.57%  '
$ws.Range('E11').Value = '  +8.85%  '
$ws.Range('D12').Value = '''0.182'
$ws.Range('E12').Value = '  +3.14%  '
$ws.Range('E13').Value = '  +8.98%  '
$ws.Range('D14').Value = '''5.43'
$ws.Range('E14').Value = '  +3.76%  '
$ws.Range('D15').Value = '''33.71'
$ws.Range('E15').Value = '  +5.06%  '
$ws.Range('D16').Value = '90.403.99'
$ws.Range('E16').Value = '  +3.56%  '
$ws.Range('D17').Value = '3.795.05'
$ws.Range('E17').Value = '  +1.29%  '
$ws.Range('D18').Value = '3.201.88'
$ws.Range('E18').Value = '  +1.51%  '
$ws.Range('B19').Value = 'PEPE'
$ws.Range('C19').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D19').Value = '''0.0000229'
$ws.Range('E19').Value = '  +78.04%  '
$ws.Range('B20').Value = 'SuiNetwork'
$ws.Range('C20').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D20').Value = '''3.37'
$ws.Range('E20').Value = '  +10.25%  '
$ws.Range('D21').Value = '''13.56'
$ws.Range('E21').Value = '  +1.23%  '
$ws.Range('D22').Value = '''441.86'
$ws.Range('E22').Value = '  +6.66%  '
$ws.Range('D23').Value = '''8.68'
$ws.Range('E23').Value = '  +3.15%  '
$ws.Range('D24').Value = '''5.09'
$ws.Range('E24').Value = '  +0.97%  '
$ws.Range('D25').Value = '''5.36'
$ws.Range('E25').Value = '  +4.12%  '
$ws.Range('D26').Value = '''11.99'
$ws.Range('E26').Value = '  +1.11%  '
$ws.Range('D27').Value = '''81.73'
$ws.Range('E27').Value = '  +11.63%  '
$ws.Range('D28').Value = '3.375.72'
$ws.Range('E28').Value = '  +1.28%  '
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('D30').Value = '''0.159'
$ws.Range('E30').Value = '  +1.77%  '
$ws.Range('D31').Value = '''0.999'
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('D32').Value = '''4.28'
$ws.Range('E32').Value = '  +43.45%  '
$ws.Range('E33').Value = '  +3.83%  '
$ws.Range('D34').Value = '''544.70'
$ws.Range('E34').Value = '  +0.36%  '
$ws.Range('D35').Value = '''7.12'
$ws.Range('E35').Value = '  +6.40%  '
$ws.Range('D36').Value = '''1.93'
$ws.Range('E36').Value = '  +4.50%  '
$ws.Range('D37').Value = '''1.31'
$ws.Range('E37').Value = '  +1.62%  '
$ws.Range('D38').Value = '''22.60'
$ws.Range('E38').Value = '  +4.00%  '
$ws.Range('E39').Value = '  +2.68%  '
$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').Value = '''1.00'
$ws.Range('E40').Value = '  +0.19%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '''0.127'
$ws.Range('E41').Value = '  -3.34%  '
$ws.Range('E42').Value = '  +3.13%  '
$ws.Range('E43').Value = '  -0.07%  '
$ws.Range('E44').Value = '  +2.28%  '
$ws.Range('D45').Value = '''146.62'
$ws.Range('E45').Value = '  -1.28%  '
$ws.Range('E46').Value = '  +4.35%  '
$ws.Range('D47').Value = '''174.17'
$ws.Range('E47').Value = '  +1.07%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').Value = '''0.758'
$ws.Range('E48').Value = '  +9.18%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').Value = '''0.126'
$ws.Range('E49').Value = '  +1.54%  '
$ws.Range('E50').Value = '  +7.66%  '
$ws.Range('D51').Value = '''1.24'
$ws.Range('E51').Value = '  +2.15%  '
